# daily auto push: 2026-01-30 02:50 UTC
#
# A new measurement was recorded for 2026/01/30 (time=8, ranking=201).
# It belongs right after the existing 2026/01/30 row (row 747), so every
# subsequent row (748..789) shifts down by one (749..790) and the sheet's
# used range grows from A1:D789 to A1:D790.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 748..789 down to 749..790, opening up a blank row 748.
$ws.Rows.Item(748).Insert()

# Column A holds a "date" value that is really free text (e.g. "2026/01/30"),
# not an Excel date serial. Force text storage, assign it, then drop the
# formatting again so the cell ends up unstyled (t="inlineStr"-equivalent,
# no number format) just like its neighbours.
$ws.Range("A748").NumberFormat = "@"
$ws.Range("A748").Value = "2026/01/30"
$ws.Range("A748").ClearFormats()

$ws.Range("B748").Value = "金"
$ws.Range("C748").Value = 8
$ws.Range("D748").Value = 201
